$p = $ppt.ActivePresentation

# The four slide titles below were stored as one run per word (split on
# spaces). Re-assigning identical-looking text is treated as a no-op by
# this host when the rendered text already matches, so nudge the value
# through an intermediate string first to force the runs to collapse
# into a single run (preserving the paragraph's pPr/buNone + rPr).
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "_"
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "Slide 1"

$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "_"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "Slide 3"

$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "_"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "Slide 4"

$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "_"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "Slide 5"

# Slide 7's speaker notes body was likewise split word-by-word; collapse
# it into a single run with the full sentence.
$notes = $p.Slides.Item(7).NotesPage.Shapes.Item(2)
$notes.TextFrame.TextRange.Text = "_"
$notes.TextFrame.TextRange.Text = "This is a blank slide: does it have a footer?"
